$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.27993369102478
$ws.Range("B1").Value = 5.660700798034668
$ws.Range("C1").Value = 2.35307240486145
$ws.Range("D1").Value = 1.542247653007507
$ws.Range("E1").Value = 1.268530130386353
